$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.801435112953186
$ws.Range("B1").Value = 1.980481386184692
$ws.Range("C1").Value = 2.194414854049683
$ws.Range("D1").Value = 3.309619188308716
$ws.Range("E1").Value = 2.045021533966064
